$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.614.98"
$ws.Range("E2").Value = "'  +3.62%  "

$ws.Range("D3").Value = "'3.248.04"
$ws.Range("E3").Value = "'  +6.45%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.11%  "

$ws.Range("D5").Value = "'581.03"
$ws.Range("E5").Value = "'  +4.41%  "

$ws.Range("D6").Value = "'153.76"
$ws.Range("E6").Value = "'  +8.41%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  -0.09%  "

$ws.Range("D8").Value = "'3.240.04"
$ws.Range("E8").Value = "'  +6.27%  "

$ws.Range("D9").Value = "'0.514"

$ws.Range("D10").Value = "'7.13"
$ws.Range("E10").Value = "'  +9.84%  "

$ws.Range("E11").Value = "'  +5.47%  "

$ws.Range("D12").Value = "'0.488"
$ws.Range("E12").Value = "'  +4.06%  "

$ws.Range("D13").Value = "'37.91"
$ws.Range("E13").Value = "'  +3.35%  "

$ws.Range("E14").Value = "'  +5.31%  "

$ws.Range("D15").Value = "'3.772.17"
$ws.Range("E15").Value = "'  +6.25%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "'66.641.24"
$ws.Range("E16").Value = "'  +3.42%  "

$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").Value = "'557.19"
$ws.Range("E17").Value = "'  +12.46%  "

$ws.Range("D18").Value = "'3.251.14"
$ws.Range("E18").Value = "'  +6.35%  "

$ws.Range("E19").Value = "'  +2.93%  "

$ws.Range("E20").Value = "'  +5.37%  "

$ws.Range("D21").Value = "'14.40"
$ws.Range("E21").Value = "'  +4.17%  "

$ws.Range("D22").Value = "'0.742"
$ws.Range("E22").Value = "'  +7.00%  "

$ws.Range("D23").Value = "'7.76"
$ws.Range("E23").Value = "'  +7.48%  "

$ws.Range("D24").Value = "'13.61"
$ws.Range("E24").Value = "'  +6.35%  "

$ws.Range("D25").Value = "'81.76"
$ws.Range("E25").Value = "'  +3.04%  "

$ws.Range("E26").Value = "'  -0.05%  "

$ws.Range("D27").Value = "'9.23"
$ws.Range("E27").Value = "'  +17.83%  "

$ws.Range("E28").Value = "'  +7.09%  "

$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "'  +5.28%  "

$ws.Range("D30").Value = "'27.76"
$ws.Range("E30").Value = "'  +6.09%  "

$ws.Range("D31").Value = "'2.76"
$ws.Range("E31").Value = "'  +4.11%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "'  -0.25%  "

$ws.Range("E33").Value = "'  +5.70%  "

$ws.Range("D34").Value = "'561.81"
$ws.Range("E34").Value = "'  +8.81%  "

$ws.Range("D35").Value = "'5.70"
$ws.Range("E35").Value = "'  +3.42%  "

$ws.Range("D36").Value = "'6.37"
$ws.Range("E36").Value = "'  +6.24%  "

$ws.Range("D37").Value = "'0.0458"
$ws.Range("E37").Value = "'  +12.82%  "

$ws.Range("D38").Value = "'55.44"
$ws.Range("E38").Value = "'  +5.35%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.132"
$ws.Range("E39").Value = "'  +7.29%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0861"
$ws.Range("E40").Value = "'  +6.73%  "

$ws.Range("D41").Value = "'3.04"
$ws.Range("E41").Value = "'  +13.33%  "

$ws.Range("D42").Value = "'3.165.63"
$ws.Range("E42").Value = "'  +8.27%  "

$ws.Range("D43").Value = "'8.61"
$ws.Range("E43").Value = "'  +1.95%  "

$ws.Range("D44").Value = "'0.274"
$ws.Range("E44").Value = "'  +10.72%  "

$ws.Range("D45").Value = "'2.31"
$ws.Range("E45").Value = "'  +8.29%  "

$ws.Range("D46").Value = "'26.51"
$ws.Range("E46").Value = "'  +4.49%  "

$ws.Range("D48").Value = "'0.0₃0556"
$ws.Range("E48").Value = "'  +3.14%  "

$ws.Range("D49").Value = "'125.98"
$ws.Range("E49").Value = "'  +3.91%  "

$ws.Range("E50").Value = "'  +2.31%  "

$ws.Range("D51").Value = "'2.23"
$ws.Range("E51").Value = "'  +7.53%  "
